# Generate Report for Handback
#
# The handback status report was regenerated: the previous run produced
# files identified by UUID "099d4bef-d7d9-445a-b17a-7624543dc0cf" (row 2,
# the canonical/non-duplicate source) and UUID
# "f15bc8f7-7e65-4ef3-9c56-d0dcd695fd4f" (row 3, flagged as a content
# duplicate of row 2). The new run produced the same two files under new
# UUIDs ("2047bb2e-acca-4a3c-8dfd-f8a610485fda" and
# "ffff3f313e35-bb77-428d-8718-d48193819a0b" respectively), with a new
# xliff correspondence hash ("fb05da26239bb986374d99f05efe0bc1d19a2bf4")
# and new handoff/handback timestamps. Because row 3 is a content
# duplicate, its "Correspond Handoff/Handback File" + datetime columns
# mirror row 2's (canonical) file, both before and after this edit.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------------
$wsOverview.Range("A2").Value = "2047bb2e-acca-4a3c-8dfd-f8a610485fda.md"
$wsOverview.Range("B2").Value = "e2e\2047bb2e-acca-4a3c-8dfd-f8a610485fda.md"
$wsOverview.Range("G2").Value = "2016-08-31 09:18:37"
$wsOverview.Range("A3").Value = "ffff3f313e35-bb77-428d-8718-d48193819a0b.md"
$wsOverview.Range("B3").Value = "e2e\ffff3f313e35-bb77-428d-8718-d48193819a0b.md"
$wsOverview.Range("G3").Value = "2016-08-31 09:18:37"

foreach ($hl in $wsOverview.Hyperlinks) {
    if ($hl.Range.Address() -eq '$B$2') { $hl.TextToDisplay = "e2e\2047bb2e-acca-4a3c-8dfd-f8a610485fda.md" }
    if ($hl.Range.Address() -eq '$B$3') { $hl.TextToDisplay = "e2e\ffff3f313e35-bb77-428d-8718-d48193819a0b.md" }
}

# --- zh-cn sheet ------------------------------------------------------------
$wsZhCn.Range("A2").Value = "2047bb2e-acca-4a3c-8dfd-f8a610485fda.md"
$wsZhCn.Range("G2").Value = "2047bb2e-acca-4a3c-8dfd-f8a610485fda.fb05da26239bb986374d99f05efe0bc1d19a2bf4.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-31 09:18:32"
$wsZhCn.Range("I2").Value = "2047bb2e-acca-4a3c-8dfd-f8a610485fda.md"
$wsZhCn.Range("J2").Value = "2047bb2e-acca-4a3c-8dfd-f8a610485fda.fb05da26239bb986374d99f05efe0bc1d19a2bf4.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-31 09:18:49"

$wsZhCn.Range("A3").Value = "ffff3f313e35-bb77-428d-8718-d48193819a0b.md"
# Row 3 is a content duplicate of row 2, so its handoff/handback file stays
# aligned to row 2's (new) xliff correspondence.
$wsZhCn.Range("G3").Value = "2047bb2e-acca-4a3c-8dfd-f8a610485fda.fb05da26239bb986374d99f05efe0bc1d19a2bf4.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-31 09:18:32"
$wsZhCn.Range("I3").Value = "ffff3f313e35-bb77-428d-8718-d48193819a0b.md"
$wsZhCn.Range("J3").Value = "2047bb2e-acca-4a3c-8dfd-f8a610485fda.fb05da26239bb986374d99f05efe0bc1d19a2bf4.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-31 09:18:49"

foreach ($hl in $wsZhCn.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') { $hl.TextToDisplay = "2047bb2e-acca-4a3c-8dfd-f8a610485fda.md" }
    if ($addr -eq '$I$2') { $hl.TextToDisplay = "2047bb2e-acca-4a3c-8dfd-f8a610485fda.md" }
    if ($addr -eq '$A$3') { $hl.TextToDisplay = "ffff3f313e35-bb77-428d-8718-d48193819a0b.md" }
    if ($addr -eq '$I$3') { $hl.TextToDisplay = "ffff3f313e35-bb77-428d-8718-d48193819a0b.md" }
}

# --- de-de sheet ------------------------------------------------------------
$wsDeDe.Range("A2").Value = "2047bb2e-acca-4a3c-8dfd-f8a610485fda.md"
$wsDeDe.Range("G2").Value = "2047bb2e-acca-4a3c-8dfd-f8a610485fda.fb05da26239bb986374d99f05efe0bc1d19a2bf4.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-31 09:18:37"
$wsDeDe.Range("I2").Value = "2047bb2e-acca-4a3c-8dfd-f8a610485fda.md"
$wsDeDe.Range("J2").Value = "2047bb2e-acca-4a3c-8dfd-f8a610485fda.fb05da26239bb986374d99f05efe0bc1d19a2bf4.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-31 09:18:56"

$wsDeDe.Range("A3").Value = "ffff3f313e35-bb77-428d-8718-d48193819a0b.md"
# Row 3 is a content duplicate of row 2, so its handoff/handback file and
# datetime stay aligned to row 2's (new) values.
$wsDeDe.Range("G3").Value = "2047bb2e-acca-4a3c-8dfd-f8a610485fda.fb05da26239bb986374d99f05efe0bc1d19a2bf4.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-31 09:18:37"
$wsDeDe.Range("I3").Value = "ffff3f313e35-bb77-428d-8718-d48193819a0b.md"
$wsDeDe.Range("J3").Value = "2047bb2e-acca-4a3c-8dfd-f8a610485fda.fb05da26239bb986374d99f05efe0bc1d19a2bf4.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-31 09:18:56"

foreach ($hl in $wsDeDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') { $hl.TextToDisplay = "2047bb2e-acca-4a3c-8dfd-f8a610485fda.md" }
    if ($addr -eq '$I$2') { $hl.TextToDisplay = "2047bb2e-acca-4a3c-8dfd-f8a610485fda.md" }
    if ($addr -eq '$A$3') { $hl.TextToDisplay = "ffff3f313e35-bb77-428d-8718-d48193819a0b.md" }
    if ($addr -eq '$I$3') { $hl.TextToDisplay = "ffff3f313e35-bb77-428d-8718-d48193819a0b.md" }
}
